# Insert a new weekly data row at row 85 (pushing existing rows 85-200 down
# to 86-201) and populate it with the new "Apio" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 85..200 down by one row, creating a blank row 85.
$ws.Rows.Item(85).Insert()

# Fill in the new row 85 with the new weekly record.
$ws.Cells.Item(85, 1).Value  = 3
$ws.Cells.Item(85, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(85, 3).Value  = 'Coquimbo'
$ws.Cells.Item(85, 4).Value  = 44413
$ws.Cells.Item(85, 5).Value  = 5
$ws.Cells.Item(85, 6).Value  = 100112017
$ws.Cells.Item(85, 7).Value  = 'Apio'
$ws.Cells.Item(85, 8).Value  = 'Americana (o)'
$ws.Cells.Item(85, 9).Value  = 'Primera'
$ws.Cells.Item(85, 10).Value = 120
$ws.Cells.Item(85, 11).Value = 9000
$ws.Cells.Item(85, 12).Value = 9000
$ws.Cells.Item(85, 13).Value = 9000
$ws.Cells.Item(85, 14).Value = '$/docena de matas'
$ws.Cells.Item(85, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(85, 16).Value = 1500
$ws.Cells.Item(85, 17).Value = 6
$ws.Cells.Item(85, 18).Value = 'Hortaliza'

# Make sure the new date cell uses the same date/time number format as the
# rest of column D.
$ws.Cells.Item(85, 4).NumberFormat = $ws.Cells.Item(86, 4).NumberFormat
